$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Context: a paragraph contains the (typo'd) tag-pair
#   "<<orr><exp>" ... "ent" ... "</exp></corr>"
# which was meant to read "<corr><exp>ent</exp></corr>". The fix
# splits the corrupted text back out into its own "<corr>"/"</corr>"
# runs (restored to the blue "<corr>"/red "</corr>" tag styling used
# everywhere else in the document) while leaving the grey "<exp>"/
# "</exp>" runs (and the untouched "ent" run) in place.
#
# Target run layout after the edit (all Courier New except "ent"):
#   1. "<corr>"   color 0000ff, size 9pt  (NEW run)
#   2. "<exp>"    color a9a9a9, size 7pt  (was "<<orr><exp>")
#   3. "ent"      (unchanged, not Courier New)
#   4. "</exp>"   color a9a9a9, size 7pt  (was "</exp></corr>")
#   5. "</corr>"  color a91111, size 9pt  (NEW run)
# ------------------------------------------------------------------

$openTag  = "<<orr><exp>"
$closeTag = "</exp></corr>"

# ==================================================================
# Opening tag: "<<orr><exp>" -> "<exp>", with a new blue "<corr>"
# run inserted immediately before it.
# ==================================================================

# Donor run carrying the exact blue/Courier-New/9pt formatting we
# need (an existing "<tl>" tag elsewhere in the document) - copying
# it (rather than hand-building a run) keeps every rPr detail
# (rFonts ascii/eastAsia/hAnsi/cs, sz, szCs, rtl) consistent with the
# rest of the file.
$donorOpen = $d.Content
$foundDonorOpen = $donorOpen.Find.Execute("<tl>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundDonorOpen) { throw "could not find donor run '<tl>'" }
$donorOpen.Copy()

# Locate the corrupted run, remember its start, and paste the donor
# run immediately before it.
$r1 = $d.Content
$found1 = $r1.Find.Execute($openTag, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "could not find '$openTag'" }
$insertPos1 = $r1.Start
$ins1 = $d.Range($insertPos1, $insertPos1)
$ins1.Paste()

# Re-point at the pasted run (same length as "<tl>") and fix its text
# to "<corr>" (formatting is left untouched).
$pastedOpen = $d.Range($insertPos1, $insertPos1 + 4)
$pastedOpen.Text = "<corr>"

# Now shrink the original corrupted text down to the plain "<exp>"
# tag it should have contained, keeping its existing formatting.
$r1fix = $d.Content
$r1fix.Find.Execute($openTag, $true, $false, $false, $false, $false, $true, 1, $false, "<exp>", 2) | Out-Null

# ==================================================================
# Closing tag: "</exp></corr>" -> "</exp>", with a new red "</corr>"
# run inserted immediately after it.
# ==================================================================

# Donor run carrying the exact red/Courier-New/9pt formatting we
# need (the pre-existing "</corr>" tag elsewhere in the document) -
# its text already matches "</corr>" so no further edit is needed
# after the paste.
$donorClose = $d.Content
$foundDonorClose = $donorClose.Find.Execute("</corr>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundDonorClose) { throw "could not find donor run '</corr>'" }
$donorClose.Copy()

# Shrink the corrupted closing text down to "</exp>" first, so the
# range's End lands exactly where the new "</corr>" run must go.
$r2 = $d.Content
$found2 = $r2.Find.Execute($closeTag, $true, $false, $false, $false, $false, $true, 1, $false, "</exp>", 2)
if (-not $found2) { throw "could not find '$closeTag'" }
$insertPos2 = $r2.End
$ins2 = $d.Range($insertPos2, $insertPos2)
$ins2.Paste()

Write-Output "done"
